# Insert a new data row before row 411 (shifts existing rows 411-459 down to 412-460)
# and populate the new row with the latest weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(411).Insert()

$ws.Cells.Item(411, 1).Value = 3
$ws.Cells.Item(411, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(411, 3).Value = "Coquimbo"
$ws.Cells.Item(411, 4).Value = 44918
$ws.Cells.Item(411, 5).Value = 5
$ws.Cells.Item(411, 6).Value = 100114013
$ws.Cells.Item(411, 7).Value = "Zanahoria"
$ws.Cells.Item(411, 8).Value = "Sin especificar"
$ws.Cells.Item(411, 9).Value = "Primera"
$ws.Cells.Item(411, 10).Value = 430
$ws.Cells.Item(411, 11).Value = 9500
$ws.Cells.Item(411, 12).Value = 10000
$ws.Cells.Item(411, 13).Value = 9756
$ws.Cells.Item(411, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(411, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(411, 16).Value = 488
$ws.Cells.Item(411, 17).Value = 20
$ws.Cells.Item(411, 18).Value = "Hortaliza"
